# Auto commit: 2024-10-14 02:06:58
# Remove the trailing period from the "DIY furniture" post's NAME (column B,
# row 4 of the "posts" sheet), and move the saved selection to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("posts")

# B4 previously held "Dive into the world of DIY furniture building using
# reclaimed wood." (with a trailing period) -- strip the trailing period.
$ws.Range("B4").Value = "Dive into the world of DIY furniture building using reclaimed wood"

# Update the sheet's saved selection/active cell to B5 (was F15).
$ws.Range("B5").Select()
